$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column headers to snake_case English names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San José De Gracia"
$ws.Range("B33").Value = "Amatenango Del Valle"
$ws.Range("B41").Value = "Comitán De Domínguez"
$ws.Range("B54").Value = "Mazapa De Madero"
$ws.Range("B57").Value = "Ocozocoautla De Espinosa"
$ws.Range("B64").Value = "Salto De Agua"
$ws.Range("B91").Value = "Coyame Del Sotol"
$ws.Range("B100").Value = "Guadalupe Y Calvo"
$ws.Range("B102").Value = "Hidalgo Del Parral"
$ws.Range("B121").Value = "San Francisco De Borja"
$ws.Range("B122").Value = "San Francisco De Conchos"
$ws.Range("B123").Value = "San Francisco Del Oro"
$ws.Range("B129").Value = "Valle De Zaragoza"
$ws.Range("B146").Value = "San Juan De Sabinas"
$ws.Range("A161").Value = "Ciudad De México"
$ws.Range("B165").Value = "Cuajimalpa De Morelos"
$ws.Range("B180").Value = "Coneto De Comonfort"
$ws.Range("B194").Value = "Nombre De Dios"
$ws.Range("B197").Value = "Pánuco De Coronado"
$ws.Range("B204").Value = "San Juan De Guadalupe"
$ws.Range("B205").Value = "San Juan Del Río"
$ws.Range("B206").Value = "San Luis Del Cordero"
$ws.Range("B207").Value = "San Pedro Del Gallo"
$ws.Range("A217").Value = "Estado De México"
$ws.Range("B217").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B218").Value = "Almoloya Del Río"
$ws.Range("B222").Value = "Atizapán De Zaragoza"
$ws.Range("B227").Value = "Coacalco De Berriozábal"
$ws.Range("B232").Value = "Ecatepec De Morelos"
$ws.Range("B236").Value = "Ixtapan De La Sal"
$ws.Range("B248").Value = "Naucalpan De Juárez"
$ws.Range("B254").Value = "San Felipe Del Progreso"
$ws.Range("B255").Value = "San Martín De Las Pirámides"
$ws.Range("B256").Value = "Soyaniquilpan De Juárez"
$ws.Range("B266").Value = "Tlalnepantla De Baz"
$ws.Range("B270").Value = "Valle De Bravo"
$ws.Range("B271").Value = "Valle De Chalco Solidaridad"
$ws.Range("B272").Value = "Villa De Allende"
$ws.Range("B282").Value = "San Miguel De Allende"
$ws.Range("B283").Value = "Apaseo El Alto"
$ws.Range("B284").Value = "Apaseo El Grande"
$ws.Range("B291").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B295").Value = "Jaral Del Progreso"
$ws.Range("B303").Value = "Purísima Del Rincón"
$ws.Range("B307").Value = "San Diego De La Unión"
$ws.Range("B309").Value = "San Francisco Del Rincón"
$ws.Range("B311").Value = "San Luis De La Paz"
$ws.Range("B312").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B313").Value = "Silao De La Victoria"
$ws.Range("B317").Value = "Valle De Santiago"
$ws.Range("B323").Value = "Acapulco De Juárez"
$ws.Range("B324").Value = "Ajuchitlán Del Progreso"
$ws.Range("B329").Value = "Atoyac De Álvarez"
$ws.Range("B330").Value = "Ayutla De Los Libres"
$ws.Range("B333").Value = "Buenavista De Cuéllar"
$ws.Range("B334").Value = "Chilapa De Álvarez"
$ws.Range("B335").Value = "Chilpancingo De Los Bravo"
$ws.Range("B338").Value = "Coyuca De Benítez"
$ws.Range("B339").Value = "Coyuca De Catalán"
$ws.Range("B343").Value = "Cutzamala De Pinzón"
$ws.Range("B349").Value = "Huitzuco De Los Figueroa"
$ws.Range("B350").Value = "Iguala De La Independencia"
$ws.Range("B352").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B353").Value = "Zihuatanejo De Azueta"
$ws.Range("B355").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B367").Value = "Taxco De Alarcón"
$ws.Range("B368").Value = "Técpan De Galeana"
$ws.Range("B371").Value = "Tixtla De Guerrero"
$ws.Range("B373").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B374").Value = "Tlapa De Comonfort"
$ws.Range("B384").Value = "Atotonilco El Grande"
$ws.Range("B391").Value = "Huasca De Ocampo"
$ws.Range("B393").Value = "Huejutla De Reyes"
$ws.Range("B396").Value = "Jacala De Ledezma"
$ws.Range("B400").Value = "Mineral Del Monte"
$ws.Range("B401").Value = "Mixquiahuala De Juárez"
$ws.Range("B402").Value = "Molango De Escamilla"
$ws.Range("B404").Value = "Nopala De Villagrán"
$ws.Range("B405").Value = "Omitlán De Juárez"
$ws.Range("B406").Value = "Pachuca De Soto"
$ws.Range("B409").Value = "Progreso De Obregón"
$ws.Range("B414").Value = "Tepehuacán De Guerrero"
$ws.Range("B415").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B416").Value = "Tezontepec De Aldama"
$ws.Range("B419").Value = "Tula De Allende"
$ws.Range("B420").Value = "Tulancingo De Bravo"
$ws.Range("B421").Value = "Zacualtipán De Ángeles"
$ws.Range("B425").Value = "Acatlán De Juárez"
$ws.Range("B426").Value = "Ahualulco De Mercado"
$ws.Range("B431").Value = "Atotonilco El Alto"
$ws.Range("B435").Value = "Cañadas De Obregón"
$ws.Range("B439").Value = "Cuautitlán De García Barragán"
$ws.Range("B444").Value = "Encarnación De Díaz"
$ws.Range("B448").Value = "Huejuquilla El Alto"
$ws.Range("B449").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B450").Value = "Ixtlahuacán Del Río"
$ws.Range("B454").Value = "Jilotlán De Los Dolores"
$ws.Range("B458").Value = "Lagos De Moreno"
$ws.Range("B462").Value = "Ojuelos De Jalisco"
$ws.Range("B466").Value = "San Cristóbal De La Barranca"
$ws.Range("B467").Value = "San Diego De Alejandría"
$ws.Range("B468").Value = "San Juan De Los Lagos"
$ws.Range("B470").Value = "San Martín De Bolaños"
$ws.Range("B471").Value = "San Miguel El Alto"
$ws.Range("B472").Value = "Santa María De Los Ángeles"
$ws.Range("B475").Value = "Talpa De Allende"
$ws.Range("B476").Value = "Tamazula De Gordiano"
$ws.Range("B478").Value = "Tepatitlán De Morelos"
$ws.Range("B479").Value = "Tlajomulco De Zúñiga"
$ws.Range("B485").Value = "Unión De San Antonio"
$ws.Range("B488").Value = "Yahualica De González Gallo"
$ws.Range("B491").Value = "Zapotitlán De Vadillo"
$ws.Range("B492").Value = "Zapotlán El Grande"
$ws.Range("B512").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B514").Value = "Cojumatlán De Régules"
$ws.Range("B562").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B581").Value = "Coatlán Del Río"
$ws.Range("B589").Value = "Puente De Ixtla"
$ws.Range("B594").Value = "Tetela Del Volcán"
$ws.Range("B601").Value = "Bahía De Banderas"
$ws.Range("B603").Value = "Ixtlán Del Río"
$ws.Range("B609").Value = "Santa María Del Oro"
$ws.Range("B624").Value = "Mier Y Noriega"
$ws.Range("B628").Value = "San Nicolás De Los Garza"
$ws.Range("B632").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B636").Value = "Ciénega De Zimatlán"
$ws.Range("B638").Value = "Cuilápam De Guerrero"
$ws.Range("B639").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B640").Value = "El Barrio De La Soledad"
$ws.Range("B641").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B642").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B643").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B645").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B649").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B650").Value = "Oaxaca De Juárez"
$ws.Range("B651").Value = "Ocotlán De Morelos"
$ws.Range("B652").Value = "Putla Villa De Guerrero"
$ws.Range("B659").Value = "San Dionisio Del Mar"
$ws.Range("B666").Value = "San José Del Progreso"
$ws.Range("B683").Value = "San Pablo Villa De Mitla"
$ws.Range("B701").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B722").Value = "Santo Domingo De Morelos"
$ws.Range("B726").Value = "Tataltepec De Valdés"
$ws.Range("B727").Value = "Teotitlán De Flores Magón"
$ws.Range("B728").Value = "Tlalixtac De Cabrera"
$ws.Range("B729").Value = "Villa De Tututepec"
$ws.Range("B731").Value = "Zimatlán De Álvarez"
$ws.Range("B741").Value = "Chalchicomula De Sesma"
$ws.Range("B747").Value = "Cuetzalan Del Progreso"
$ws.Range("B756").Value = "Ixcamilpa De Guerrero"
$ws.Range("B759").Value = "Izúcar De Matamoros"
$ws.Range("B765").Value = "Palmar De Bravo"
$ws.Range("B773").Value = "San Salvador El Seco"
$ws.Range("B775").Value = "Tecali De Herrera"
$ws.Range("B779").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B781").Value = "Tepexi De Rodríguez"
$ws.Range("B784").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B792").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B797").Value = "Amealco De Bonfil"
$ws.Range("B799").Value = "Cadereyta De Montes"
$ws.Range("B804").Value = "Jalpan De Serra"
$ws.Range("B805").Value = "Landa De Matamoros"
$ws.Range("B806").Value = "Pinal De Amoles"
$ws.Range("B808").Value = "San Juan Del Río"
$ws.Range("B819").Value = "Armadillo De Los Infante"
$ws.Range("B825").Value = "Ciudad Del Maíz"
$ws.Range("B832").Value = "Mexquitic De Carmona"
$ws.Range("B837").Value = "San Ciro De Acosta"
$ws.Range("B842").Value = "Santa María Del Río"
$ws.Range("B844").Value = "Soledad De Graciano Sánchez"
$ws.Range("B851").Value = "Villa De Arista"
$ws.Range("B852").Value = "Villa De Arriaga"
$ws.Range("B853").Value = "Villa De Guadalupe"
$ws.Range("B854").Value = "Villa De La Paz"
$ws.Range("B855").Value = "Villa De Ramos"
$ws.Range("B856").Value = "Villa De Reyes"
$ws.Range("B888").Value = "Nacozari De García"
$ws.Range("B902").Value = "Jalpa De Méndez"
$ws.Range("B931").Value = "Soto La Marina"
$ws.Range("B939").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B940").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B944").Value = "San Pablo Del Monte"
$ws.Range("B954").Value = "Amatlán De Los Reyes"
$ws.Range("B960").Value = "Boca Del Río"
$ws.Range("B962").Value = "Castillo De Teayo"
$ws.Range("B964").Value = "Cazones De Herrera"
$ws.Range("B967").Value = "Chinampa De Gorostiza"
$ws.Range("B975").Value = "Cosamaloapan De Carpio"
$ws.Range("B983").Value = "Hueyapan De Ocampo"
$ws.Range("B984").Value = "Ignacio De La Llave"
$ws.Range("B986").Value = "Ixhuatlán De Madero"
$ws.Range("B987").Value = "Ixhuatlán Del Sureste"
$ws.Range("B996").Value = "Lerdo De Tejada"
$ws.Range("B997").Value = "Martínez De La Torre"
$ws.Range("B999").Value = "Medellín De Bravo"
$ws.Range("B1002").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1007").Value = "Ozuluama De Mascareñas"
$ws.Range("B1010").Value = "Paso De Ovejas"
$ws.Range("B1013").Value = "Poza Rica De Hidalgo"
$ws.Range("B1019").Value = "Sayula De Alemán"
$ws.Range("B1026").Value = "Tatahuicapan De Juárez"
$ws.Range("B1040").Value = "Vega De Alatorre"
$ws.Range("B1045").Value = "Zozocolco De Hidalgo"
$ws.Range("B1057").Value = "Cañitas De Felipe Pescador"
$ws.Range("B1059").Value = "Concepción Del Oro"
$ws.Range("B1069").Value = "Jiménez Del Teul"
$ws.Range("B1080").Value = "Moyahua De Estrada"
$ws.Range("B1081").Value = "Nochistlán De Mejía"
$ws.Range("B1082").Value = "Noria De Ángeles"
$ws.Range("B1092").Value = "Teúl De González Ortega"
$ws.Range("B1093").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1096").Value = "Villa De Cos"

# Remove the trailing metadata/footnote rows (1105-1109); row 1104 is already blank.
$ws.Range("A1104:A1109").EntireRow.Delete()
